# The deck ships two theme parts:
#   ppt/theme/theme2.xml -> referenced by the slide master / presentation (was "Integral" / Red Violet)
#   ppt/theme/theme1.xml -> referenced by the notes master (was "Office Theme")
# The target edit swaps the color schemes carried by those two theme parts:
#   the slide master's theme becomes the stock "Office Theme" colors
#   the notes master's theme becomes the "Integral" / Red Violet colors
# Font scheme / format scheme are identical between the two themes, so only
# the 10 theme colors that differ (dk2, lt2, accent1-6, hlink, folHlink) need
# to move; dk1/lt1 are already 000000/FFFFFF in both.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Slide / slide-master theme (backs ppt/theme/theme2.xml): Integral -> Office Theme
$slideColors = $s.ThemeColorScheme
$slideColors.Item(3).RGB  = 6968388    # dk2      44546A
$slideColors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$slideColors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$slideColors.Item(6).RGB  = 3243501    # accent2  ED7D31
$slideColors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$slideColors.Item(8).RGB  = 49407      # accent4  FFC000
$slideColors.Item(9).RGB  = 12874308   # accent5  4472C4
$slideColors.Item(10).RGB = 4697456    # accent6  70AD47
$slideColors.Item(11).RGB = 12673797   # hlink    0563C1
$slideColors.Item(12).RGB = 7491477    # folHlink 954F72

# --- Notes-master theme (backs ppt/theme/theme1.xml): Office Theme -> Integral
$notesColors = $s.NotesPage.ThemeColorScheme
$notesColors.Item(3).RGB  = 5326149    # dk2      454551
$notesColors.Item(4).RGB  = 14473688   # lt2      D8D9DC
$notesColors.Item(5).RGB  = 9514467    # accent1  E32D91
$notesColors.Item(6).RGB  = 13381832   # accent2  C830CC
$notesColors.Item(7).RGB  = 14460494   # accent3  4EA6DC
$notesColors.Item(8).RGB  = 15168839   # accent4  4775E7
$notesColors.Item(9).RGB  = 14774665   # accent5  8971E1
$notesColors.Item(10).RGB = 7555029    # accent6  D54773
$notesColors.Item(11).RGB = 2465643    # hlink    6B9F25
$notesColors.Item(12).RGB = 9211020    # folHlink 8C8C8C
